# Generate Report for Archive
# - Update status text "Ready for handoff" -> "In Translation" on all three sheets
# - Shrink the "Status" column(s) width to match the freshly generated report layout

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns (E and F) ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"

# Column width: set as close as possible to the generated report's target width.
$ws1.Range("E1").ColumnWidth = 12.5
$ws1.Range("F1").ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C) ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"
$ws2.Range("C3").Value = "In Translation"
$ws2.Range("C1").ColumnWidth = 12.5

# --- de-de sheet: Status column (C) ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"
$ws3.Range("C3").Value = "In Translation"
$ws3.Range("C1").ColumnWidth = 12.5
